$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 795.1781
$ws.Range("I17").Value = 964.5
$ws.Range("J17").Value = 790.40845
$ws.Range("K17").Value = 2893.5
$ws.Range("L17").Value = 2371.22535
$ws.Range("M17").Value = -2725.5
$ws.Range("N17").Value = -2707.22535
$ws.Range("H62").Value = 16885.895
$ws.Range("I62").Value = 12926.625
$ws.Range("J62").Value = 38002
$ws.Range("K62").Value = 12926.625
$ws.Range("L62").Value = 38002
$ws.Range("M62").Value = -12302.625
$ws.Range("N62").Value = -39250
$ws.Range("H65").Value = 16885.895
$ws.Range("I65").Value = 12926.625
$ws.Range("J65").Value = 38002
$ws.Range("K65").Value = 64633.125
$ws.Range("L65").Value = 190010
$ws.Range("M65").Value = -61513.125
$ws.Range("N65").Value = -196250
$ws.Range("H80").Value = 6945768
$ws.Range("J80").Value = 1799
$ws.Range("L80").Value = 5397
$ws.Range("N80").Value = -7393
$ws.Range("H83").Value = 6945768
$ws.Range("J83").Value = 1799
$ws.Range("L83").Value = 16191
$ws.Range("N83").Value = -26175
$ws.Range("H96").Value = 2044.3334
$ws.Range("J96").Value = 3509.6667
$ws.Range("L96").Value = 10529.0001
$ws.Range("N96").Value = -13275.0001
$ws.Range("H98").Value = 1071.5714
$ws.Range("I98").Value = 1000.1539
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1000.1539
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 497.8461
$ws.Range("N98").Value = -4996
$ws.Range("H99").Value = 1532
$ws.Range("I99").Value = 381.5
$ws.Range("K99").Value = 1144.5
$ws.Range("M99").Value = 353.5
$ws.Range("H111").Value = 3982.2727
$ws.Range("I111").Value = 3282.7144
$ws.Range("J111").Value = 5206.5
$ws.Range("K111").Value = 9848.143199999999
$ws.Range("L111").Value = 15619.5
$ws.Range("M111").Value = -6781.143199999999
$ws.Range("N111").Value = -21753.5
$ws.Range("H113").Value = 3581.2727
$ws.Range("J113").Value = 4266.6665
$ws.Range("L113").Value = 4266.6665
$ws.Range("N113").Value = -10774.6665
$ws.Range("H116").Value = 150633.16
$ws.Range("J116").Value = 186428.95
$ws.Range("L116").Value = 186428.95
$ws.Range("N116").Value = -193312.95
$ws.Range("H121").Value = 9599.6
$ws.Range("I121").Value = 1000
$ws.Range("K121").Value = 3000
$ws.Range("M121").Value = -1253
$ws.Range("H122").Value = 1071.5714
$ws.Range("I122").Value = 1000.1539
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000.4617
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550.4616999999998
$ws.Range("N122").Value = -10900
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6579767
$ws.Range("I32").Value = 7042935
$ws.Range("K32").Value = 7042935
$ws.Range("M32").Value = -7042648
$ws.Range("H45").Value = 2056.7334
$ws.Range("I45").Value = 2046.5
$ws.Range("K45").Value = 2046.5
$ws.Range("M45").Value = -1669.5
$ws.Range("H61").Value = 815824.3
$ws.Range("I61").Value = 1044046.75
$ws.Range("J61").Value = 4367
$ws.Range("K61").Value = 1044046.75
$ws.Range("L61").Value = 4367
$ws.Range("M61").Value = -1043834.75
$ws.Range("N61").Value = -4791
$ws.Range("H102").Value = 30766.166
$ws.Range("I102").Value = 33290.453
$ws.Range("K102").Value = 33290.453
$ws.Range("M102").Value = -31668.453
$ws.Range("H122").Value = 3642
$ws.Range("I122").Value = 3581.6086
$ws.Range("K122").Value = 10744.8258
$ws.Range("M122").Value = -8294.825800000001
$ws.Range("H132").Value = 366728
$ws.Range("I132").Value = 403201.38
$ws.Range("K132").Value = 1209604.14
$ws.Range("M132").Value = -1207074.14
$ws.Range("H136").Value = 815824.3
$ws.Range("I136").Value = 1044046.75
$ws.Range("J136").Value = 4367
$ws.Range("K136").Value = 3132140.25
$ws.Range("L136").Value = 13101
$ws.Range("M136").Value = -3129590.25
$ws.Range("N136").Value = -18201
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2436.6924
$ws.Range("J80").Value = 3226.3333
$ws.Range("L80").Value = 3226.3333
$ws.Range("N80").Value = -5222.3333
$ws.Range("H83").Value = 2436.6924
$ws.Range("J83").Value = 3226.3333
$ws.Range("L83").Value = 16131.6665
$ws.Range("N83").Value = -26115.6665
$ws.Range("H99").Value = 12956.286
$ws.Range("I99").Value = 6125.1816
$ws.Range("K99").Value = 6125.1816
$ws.Range("M99").Value = -4627.1816
$ws.Range("H107").Value = 2243.1614
$ws.Range("I107").Value = 2104.3103
$ws.Range("K107").Value = 2104.3103
$ws.Range("M107").Value = -184.3103000000001
$ws.Range("H134").Value = 590118.75
$ws.Range("I134").Value = 806946.25
$ws.Range("J134").Value = 247759.58
$ws.Range("K134").Value = 2420838.75
$ws.Range("L134").Value = 743278.74
$ws.Range("M134").Value = -2418303.75
$ws.Range("N134").Value = -748348.74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2971.111
$ws.Range("I122").Value = 1737
$ws.Range("J122").Value = 3958.4
$ws.Range("K122").Value = 5211
$ws.Range("L122").Value = 11875.2
$ws.Range("M122").Value = -2761
$ws.Range("N122").Value = -16775.2
$ws.Range("H132").Value = 22419954
$ws.Range("I132").Value = 26327234
$ws.Range("K132").Value = 78981702
$ws.Range("M132").Value = -78979172
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 16941.25
$ws.Range("I3").Value = 12510
$ws.Range("J3").Value = 19600
$ws.Range("K3").Value = 37530
$ws.Range("L3").Value = 58800
$ws.Range("M3").Value = -37418
$ws.Range("N3").Value = -59024
$ws.Range("H37").Value = 96082
$ws.Range("J37").Value = 96082
$ws.Range("L37").Value = 288246
$ws.Range("N37").Value = -288470
$ws.Range("H113").Value = 5973
$ws.Range("I113").Value = 459.5
$ws.Range("J113").Value = 17000
$ws.Range("K113").Value = 1378.5
$ws.Range("L113").Value = 51000
$ws.Range("M113").Value = 791.5
$ws.Range("N113").Value = -55340
$ws.Range("H134").Value = 975
$ws.Range("I134").Value = 975
$ws.Range("K134").Value = 2925
$ws.Range("M134").Value = 2145
$ws.Range("H141").Value = 2962.8333
$ws.Range("I141").Value = 2962.8333
$ws.Range("K141").Value = 8888.499899999999
$ws.Range("M141").Value = -3708.499899999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 169031.27
$ws.Range("I80").Value = 280114.06
$ws.Range("J80").Value = 2407.0833
$ws.Range("K80").Value = 280114.06
$ws.Range("L80").Value = 2407.0833
$ws.Range("M80").Value = -279116.06
$ws.Range("N80").Value = -4403.0833
$ws.Range("H83").Value = 169031.27
$ws.Range("I83").Value = 280114.06
$ws.Range("J83").Value = 2407.0833
$ws.Range("K83").Value = 1400570.3
$ws.Range("L83").Value = 12035.4165
$ws.Range("M83").Value = -1395578.3
$ws.Range("N83").Value = -22019.4165
$ws.Range("H126").Value = 879994.9399999999
$ws.Range("I126").Value = 1668487.5
$ws.Range("J126").Value = 3892.111
$ws.Range("K126").Value = 5005462.5
$ws.Range("L126").Value = 11676.333
$ws.Range("M126").Value = -5002992.5
$ws.Range("N126").Value = -16616.333
$ws.Range("H132").Value = 191725.72
$ws.Range("I132").Value = 267943.16
$ws.Range("K132").Value = 803829.48
$ws.Range("M132").Value = -801299.48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1949.875
$ws.Range("I68").Value = 1780
$ws.Range("J68").Value = 2233
$ws.Range("K68").Value = 1780
$ws.Range("L68").Value = 2233
$ws.Range("M68").Value = -1031
$ws.Range("N68").Value = -3731
$ws.Range("H71").Value = 1949.875
$ws.Range("I71").Value = 1780
$ws.Range("J71").Value = 2233
$ws.Range("K71").Value = 8900
$ws.Range("L71").Value = 11165
$ws.Range("M71").Value = -5156
$ws.Range("N71").Value = -18653
$ws.Range("H137").Value = 110000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 35916.668
$ws.Range("J5").Value = 41555.555
$ws.Range("L5").Value = 41555.555
$ws.Range("N5").Value = -41779.555
$ws.Range("H122").Value = 2604.1277
$ws.Range("I122").Value = 2276.5715
$ws.Range("J122").Value = 5355.6
$ws.Range("K122").Value = 6829.7145
$ws.Range("L122").Value = 16066.8
$ws.Range("M122").Value = -4379.7145
$ws.Range("N122").Value = -20966.8
$ws.Range("H126").Value = 2242.375
$ws.Range("I126").Value = 1163.3334
$ws.Range("J126").Value = 5479.5
$ws.Range("K126").Value = 3490.0002
$ws.Range("L126").Value = 16438.5
$ws.Range("M126").Value = -1020.0002
$ws.Range("N126").Value = -21378.5
$ws.Range("H132").Value = 13839297
$ws.Range("I132").Value = 21119616
$ws.Range("J132").Value = 6690.9
$ws.Range("K132").Value = 63358848
$ws.Range("L132").Value = 20072.7
$ws.Range("M132").Value = -63356318
